# Auto-generated Excel COM-interop script
# Applies market-price / profit recalculation updates to the Seraph_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow
$ws.Range("H12").Value = 802.46155
$ws.Range("I12").Value = 720.4286
$ws.Range("J12").Value = 898.1667
$ws.Range("K12").Value = 720.4286
$ws.Range("L12").Value = 898.1667
$ws.Range("M12").Value = -550.4286
$ws.Range("N12").Value = -1238.1667

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 6121.9473
$ws.Range("I76").Value = 4814.1113
$ws.Range("K76").Value = 4814.1113
$ws.Range("M76").Value = -4499.1113

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 6121.9473
$ws.Range("I79").Value = 4814.1113
$ws.Range("K79").Value = 4814.1113
$ws.Range("M79").Value = -3722.1113

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 2652.7144
$ws.Range("J111").Value = 1600
$ws.Range("L111").Value = 4800
$ws.Range("N111").Value = -10934

# Row 116: Growing Up
$ws.Range("H116").Value = 8190.8887
$ws.Range("J116").Value = 7974
$ws.Range("L116").Value = 7974
$ws.Range("N116").Value = -14858

# Row 138: All-night Crafting
$ws.Range("H138").Value = 4041.6604
$ws.Range("J138").Value = 4245.1226
$ws.Range("L138").Value = 12735.3678
$ws.Range("N138").Value = -23015.3678

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 16659.277
$ws.Range("I32").Value = 14843.615
$ws.Range("K32").Value = 14843.615
$ws.Range("M32").Value = -14556.615

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1975.8889
$ws.Range("J45").Value = 2031.6666
$ws.Range("L45").Value = 2031.6666
$ws.Range("N45").Value = -2785.6666

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 1895
$ws.Range("I61").Value = 1895
$ws.Range("K61").Value = 1895
$ws.Range("M61").Value = -1683

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 4688.0415
$ws.Range("I63").Value = 3567.6667
$ws.Range("K63").Value = 3567.6667
$ws.Range("M63").Value = -2881.6667

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 4688.0415
$ws.Range("I66").Value = 3567.6667
$ws.Range("K66").Value = 17838.3335
$ws.Range("M66").Value = -14406.3335

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 626947.6
$ws.Range("I122").Value = 715833
$ws.Range("K122").Value = 2147499
$ws.Range("M122").Value = -2145049

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 1895
$ws.Range("I136").Value = 1895
$ws.Range("K136").Value = 5685
$ws.Range("M136").Value = -3135

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 2233.0667
$ws.Range("J94").Value = 1825
$ws.Range("L94").Value = 1825
$ws.Range("N94").Value = -2727

# Row 107: The Gold Experience
$ws.Range("H107").Value = 950.625
$ws.Range("I107").Value = 950.625
$ws.Range("K107").Value = 950.625
$ws.Range("M107").Value = 969.375

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 2048.4285
$ws.Range("I16").Value = 1334.75
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1334.75
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1047.75
$ws.Range("N16").Value = -3574

# Row 31: Wall Not Found
$ws.Range("H31").Value = 4403.3687
$ws.Range("I31").Value = 1860
$ws.Range("K31").Value = 1860
$ws.Range("M31").Value = -1565

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 4403.3687
$ws.Range("I34").Value = 1860
$ws.Range("K34").Value = 1860
$ws.Range("M34").Value = -1658

# Row 107: Built to Last
$ws.Range("H107").Value = 990.9286
$ws.Range("I107").Value = 790.6667
$ws.Range("K107").Value = 790.6667
$ws.Range("M107").Value = 1129.3333

# Row 113: Patient Patients
$ws.Range("H113").Value = 2048.4285
$ws.Range("I113").Value = 1334.75
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1334.75
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 835.25
$ws.Range("N113").Value = -7340

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 20759.363
$ws.Range("J141").Value = 20759.363
$ws.Range("L141").Value = 20759.363
$ws.Range("N141").Value = -31119.363

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 28.875
$ws.Range("J12").Value = 30.333334
$ws.Range("L12").Value = 91.00000199999999
$ws.Range("N12").Value = -437.000002

# Row 122: Salt of the North
$ws.Range("H122").Value = 655.6923
$ws.Range("I122").Value = 608.6
$ws.Range("J122").Value = 685.125
$ws.Range("K122").Value = 5477.400000000001
$ws.Range("L122").Value = 6166.125
$ws.Range("M122").Value = -3027.400000000001
$ws.Range("N122").Value = -11066.125

# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 2569.7334
$ws.Range("I140").Value = 2569.7334
$ws.Range("K140").Value = 7709.2002
$ws.Range("M140").Value = -2529.2002

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 14998.667
$ws.Range("I80").Value = 14998
$ws.Range("J80").Value = 14999
$ws.Range("K80").Value = 14998
$ws.Range("L80").Value = 14999
$ws.Range("M80").Value = -14000
$ws.Range("N80").Value = -16995

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 14998.667
$ws.Range("I83").Value = 14998
$ws.Range("J83").Value = 14999
$ws.Range("K83").Value = 74990
$ws.Range("L83").Value = 74995
$ws.Range("M83").Value = -69998
$ws.Range("N83").Value = -84979

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 4601.1665
$ws.Range("I113").Value = 3879.5
$ws.Range("K113").Value = 3879.5
$ws.Range("M113").Value = -1709.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 1921.625
$ws.Range("I132").Value = 1507.3572
$ws.Range("K132").Value = 4522.071599999999
$ws.Range("M132").Value = -1992.071599999999

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 10380.777
$ws.Range("I61").Value = 9678.625
$ws.Range("J61").Value = 15998
$ws.Range("K61").Value = 9678.625
$ws.Range("L61").Value = 15998
$ws.Range("M61").Value = -9476.625
$ws.Range("N61").Value = -16402

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 3762.625
$ws.Range("I100").Value = 2266.889
$ws.Range("J100").Value = 5685.7144
$ws.Range("K100").Value = 2266.889
$ws.Range("L100").Value = 5685.7144
$ws.Range("M100").Value = -1725.889
$ws.Range("N100").Value = -6767.7144

# Row 113: Peace in Rest
$ws.Range("H113").Value = 10380.777
$ws.Range("I113").Value = 9678.625
$ws.Range("J113").Value = 15998
$ws.Range("K113").Value = 9678.625
$ws.Range("L113").Value = 15998
$ws.Range("M113").Value = -7508.625
$ws.Range("N113").Value = -20338

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3692.3462
$ws.Range("I132").Value = 3477.318
$ws.Range("J132").Value = 4875
$ws.Range("K132").Value = 10431.954
$ws.Range("L132").Value = 14625
$ws.Range("M132").Value = -7901.954000000002
$ws.Range("N132").Value = -19685

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 3901.3333
$ws.Range("I136").Value = 3901.3333
$ws.Range("K136").Value = 11703.9999
$ws.Range("M136").Value = -9153.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 47: The Wages of Sin
$ws.Range("H47").Value = 39989.668
$ws.Range("J47").Value = 39989.668
$ws.Range("L47").Value = 39989.668
$ws.Range("N47").Value = -41133.668

# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 7625
$ws.Range("J62").Value = 7708.778
$ws.Range("L62").Value = 7708.778
$ws.Range("N62").Value = -8956.778

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 7625
$ws.Range("J65").Value = 7708.778
$ws.Range("L65").Value = 38543.89
$ws.Range("N65").Value = -44783.89

# Row 95: Duress Rehearsal
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492

# Row 107: Flax Wax
$ws.Range("H107").Value = 1063.9
$ws.Range("I107").Value = 435.5
$ws.Range("K107").Value = 1306.5
$ws.Range("M107").Value = 613.5

# Row 113: A Tender Table
$ws.Range("H113").Value = 878.9375
$ws.Range("I113").Value = 743.5454999999999
$ws.Range("K113").Value = 2230.6365
$ws.Range("M113").Value = -60.63649999999961

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2635
$ws.Range("I126").Value = 805
$ws.Range("J126").Value = 8125
$ws.Range("K126").Value = 2415
$ws.Range("L126").Value = 24375
$ws.Range("M126").Value = 55
$ws.Range("N126").Value = -29315

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 3042
$ws.Range("I132").Value = 2389.3333
$ws.Range("K132").Value = 7167.999899999999
$ws.Range("M132").Value = -4637.999899999999
